# "UI Validations added V.48"
#
# The stale Viax-order test data that had been left behind in the
# UIValidations sheet (rows 3 & 4, columns B:D) is removed - those were
# left-over order numbers / "Created Viax order ..." messages from a
# previous test run that no longer belong in the validation sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UIValidations")

# Row 3 ("OrderID"): drop the leftover order numbers in B3:D3.
$ws.Range("B3:D3").ClearContents()

# Row 4 ("OrderStatus"): drop the leftover Viax order messages in B4:D4.
$ws.Range("B4:D4").ClearContents()

# Column C keeps an (empty) cell in both rows instead of being fully
# blank - touch its style so the cell stays on the sheet.
$ws.Range("C3").Style = "Normal"
$ws.Range("C4").Style = "Normal"

# Restore the view / selection to where the author left it last.
$excel.Goto($ws.Range("A43"), $true)
$ws.Range("D17").Select()
